# Chiffres COVID-19 Valais - data upload update
#
# Fills in newly-reported daily figures (rows 611, 623-632) on the
# "Feuil1" sheet and updates the active-cell selection left behind by
# the editor. Columns B, H, J and K are live cached formulas
# (IF(TODAY()>..., ...)) and recompute automatically from the literal
# inputs below, so only the literal data columns need to be written.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns L/M ("Nb nouveaux deces a l'hopital" / "...extra-hospitaliers")
# are formatted as Text (numFmtId 49) in this sheet, so a plain
# Range.Value2 write of a number gets stored as a text string. The
# existing data in the sheet stores these as real numbers, so nudge the
# number format to a plain numeric one, write the value, then restore
# the original (bordered) number format - this keeps the cell's
# original look while writing a genuine numeric value, same as the
# surrounding cells.
function Set-NumericValue($range, $val) {
    $origFormat = $range.NumberFormat
    $range.NumberFormat = "0"
    $range.Value2 = $val
    $range.NumberFormat = $origFormat
}

# Row 611 - corrected new-case count
$ws.Range("C611").Value2 = 60

# Row 623 - fill in previously-blank daily death columns with 0
Set-NumericValue $ws.Range("L623") 0
Set-NumericValue $ws.Range("M623") 0

# Row 624 - fill in previously-blank daily death columns with 0
Set-NumericValue $ws.Range("L624") 0
Set-NumericValue $ws.Range("M624") 0

# Row 625 - updated SI occupancy figures + fill blank death columns
$ws.Range("G625").Value2 = 15
Set-NumericValue $ws.Range("L625") 0
Set-NumericValue $ws.Range("M625") 0

# Row 626 - updated case/hospital counts + fill blank death columns
$ws.Range("C626").Value2 = 131
$ws.Range("E626").Value2 = 4
$ws.Range("G626").Value2 = 17
Set-NumericValue $ws.Range("L626") 0
Set-NumericValue $ws.Range("M626") 0

# Row 627 - updated case/hospital counts + fill blank death columns
$ws.Range("C627").Value2 = 77
$ws.Range("E627").Value2 = 4
$ws.Range("G627").Value2 = 17
Set-NumericValue $ws.Range("L627") 0
Set-NumericValue $ws.Range("M627") 0

# Row 628 - updated case/hospital counts + fill blank death columns
$ws.Range("C628").Value2 = 52
$ws.Range("E628").Value2 = 4
$ws.Range("G628").Value2 = 19
Set-NumericValue $ws.Range("L628") 0
Set-NumericValue $ws.Range("M628") 0

# Row 629 - updated case/hospital counts + fill blank death columns
$ws.Range("C629").Value2 = 155
$ws.Range("E629").Value2 = 4
$ws.Range("G629").Value2 = 20
Set-NumericValue $ws.Range("L629") 0
Set-NumericValue $ws.Range("M629") 0

# Row 630 - newly reported day (was completely blank before)
$ws.Range("C630").Value2 = 139
$ws.Range("E630").Value2 = 5
$ws.Range("F630").Value2 = 2
$ws.Range("G630").Value2 = 23
Set-NumericValue $ws.Range("L630") 1
Set-NumericValue $ws.Range("M630") 1

# Row 631 - newly reported day (was completely blank before)
$ws.Range("C631").Value2 = 132
$ws.Range("E631").Value2 = 6
$ws.Range("F631").Value2 = 3
$ws.Range("G631").Value2 = 29
Set-NumericValue $ws.Range("L631") 0
Set-NumericValue $ws.Range("M631") 0

# Row 632 - newly reported day (was completely blank before)
$ws.Range("C632").Value2 = 8
$ws.Range("E632").Value2 = 5
$ws.Range("F632").Value2 = 3
$ws.Range("G632").Value2 = 34
Set-NumericValue $ws.Range("L632") 0
Set-NumericValue $ws.Range("M632") 0

# Leave behind the same active-cell selection the editor had when they
# finished working (bottom-right frozen pane, cell Q620).
$ws.Range("Q620").Select()
